# "viz y tablas update"
#
# On the "Ficha técnica" sheet:
#   - remove the "DIMENSIÓN" / "Accesibilidad" row (row 3), shifting the
#     remaining metadata rows up by one
#   - append two new metadata rows at the bottom:
#       TIPOIND / Resultados
#       CITA    / UMAD con base en Instituto de Economía, Universidad de la
#                 República (2020) Encuesta Continua de Hogares
#                 Compatibilizada 1981-2018 Versión 12 DOI:
#                 http://doiorg/1047426/ECHINE

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ficha técnica")

# Drop the "DIMENSIÓN" / "Accesibilidad" row entirely; rows below shift up.
$ws.Rows(3).Delete()

# New metadata rows appended after "CÁLCULO" (now row 6).
$ws.Range("A7").Value = "TIPOIND"
$ws.Range("B7").Value = "Resultados"
$ws.Range("A8").Value = "CITA"
$ws.Range("B8").Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"
